$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Resolving-Mac" target-cluster row (was row 5) entirely.
$ws.Rows.Item(5).Delete()

# Row 2 (target cluster: ECs) - updated TPM-derived values
$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("M2").Value = 30.58864766666666
$ws.Range("N2").Value = 91.76594299999999
$ws.Range("O2").Value = 0.3925391465174898
$ws.Range("P2").Value = 0.3925391465174898
$ws.Range("Q2").Value = 20.83469264195833
$ws.Range("R2").Value = 187.512233777625
$ws.Range("S2").Value = 0.3925391465174898
$ws.Range("T2").Value = 0.3925391465174898

# Row 3 (target cluster: FAPs) - updated TPM-derived values
$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.291183949679193
$ws.Range("P3").Value = 0.291183949679193
$ws.Range("S3").Value = 0.291183949679193
$ws.Range("T3").Value = 0.291183949679193

# Row 4 (target cluster: MuSCs) - updated TPM-derived values
$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 24.64590566666666
$ws.Range("N4").Value = 73.93771699999999
$ws.Range("O4").Value = 0.3162769038033173
$ws.Range("P4").Value = 0.3162769038033172
$ws.Range("Q4").Value = 16.78694249720833
$ws.Range("R4").Value = 151.082482474875
$ws.Range("S4").Value = 0.3162769038033173
$ws.Range("T4").Value = 0.3162769038033172
